$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data point was appended to this rolling fuel-price time series.
# Row 3 (the blank placeholder row that used to sit above the most recent
# priced row) is removed, shifting every row below it up by one position,
# and the row that becomes the new row 4 is populated with the newly
# published price figures.
$ws.Rows("3").Delete()

$ws.Range("A4").Value = 770.419
$ws.Range("B4").Value = 717.028
